$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.871.25'
$ws.Range('E2').Value = '  +0.74%  '
$ws.Range('D3').Value = '3.834.69'
$ws.Range('E3').Value = '  -0.60%  '
$ws.Range('E4').Value = '  +0.02%  '
$r = $ws.Range('D5')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '704.45'
$r.Style = $origStyle
$ws.Range('E5').Value = '  -0.01%  '
$r = $ws.Range('D6')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '172.40'
$r.Style = $origStyle
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('D7').Value = '3.836.23'
$ws.Range('E7').Value = '  -0.49%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('E10').Value = '  -2.18%  '
$r = $ws.Range('D11')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '7.47'
$r.Style = $origStyle
$ws.Range('E11').Value = '  +2.28%  '
$r = $ws.Range('D12')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.486'
$r.Style = $origStyle
$ws.Range('E12').Value = '  +5.59%  '
$r = $ws.Range('D13')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.0000254'
$r.Style = $origStyle
$ws.Range('E13').Value = '  -2.01%  '
$r = $ws.Range('D14')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '36.54'
$r.Style = $origStyle
$ws.Range('E14').Value = '  -0.72%  '
$ws.Range('D15').Value = '4.483.81'
$ws.Range('E15').Value = '  -0.42%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '71.922.92'
$ws.Range('E16').Value = '  +0.92%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.793.73'
$ws.Range('E17').Value = '  -0.04%  '
$r = $ws.Range('D18')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '7.25'
$r.Style = $origStyle
$ws.Range('E18').Value = '  +0.11%  '
$r = $ws.Range('D19')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '17.58'
$r.Style = $origStyle
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('E20').Value = '  -0.21%  '
$r = $ws.Range('D21')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '518.71'
$r.Style = $origStyle
$ws.Range('E21').Value = '  +4.38%  '
$r = $ws.Range('D22')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '10.61'
$r.Style = $origStyle
$ws.Range('E22').Value = '  -1.04%  '
$r = $ws.Range('D23')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.720'
$r.Style = $origStyle
$ws.Range('E23').Value = '  -1.40%  '
$r = $ws.Range('D24')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '84.51'
$r.Style = $origStyle
$ws.Range('E24').Value = '  -1.04%  '
$r = $ws.Range('D25')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.0000143'
$r.Style = $origStyle
$ws.Range('E25').Value = '  -3.58%  '
$r = $ws.Range('D26')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '12.75'
$r.Style = $origStyle
$ws.Range('E26').Value = '  +4.17%  '
$ws.Range('D27').Value = '3.984.52'
$ws.Range('E27').Value = '  -0.56%  '
$r = $ws.Range('D28')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '10.45'
$r.Style = $origStyle
$ws.Range('E28').Value = '  -2.36%  '
$ws.Range('E29').Value = '  +0.02%  '
$r = $ws.Range('D30')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.04'
$r.Style = $origStyle
$ws.Range('E30').Value = '  -3.59%  '
$r = $ws.Range('D31')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '3.04'
$r.Style = $origStyle
$ws.Range('E31').Value = '  -5.05%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$r = $ws.Range('D32')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.26'
$r.Style = $origStyle
$ws.Range('E32').Value = '  -0.58%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$r = $ws.Range('D33')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '7.41'
$r.Style = $origStyle
$ws.Range('E33').Value = '  -1.84%  '
$r = $ws.Range('D34')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '29.37'
$r.Style = $origStyle
$ws.Range('E34').Value = '  -0.74%  '
$r = $ws.Range('D35')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.175'
$r.Style = $origStyle
$ws.Range('E35').Value = '  -2.76%  '
$r = $ws.Range('D36')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '9.25'
$r.Style = $origStyle
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '3.800.21'
$ws.Range('E37').Value = '  -0.34%  '
$r = $ws.Range('D38')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.997'
$r.Style = $origStyle
$ws.Range('E38').Value = '  -0.28%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$r = $ws.Range('D39')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '2.49'
$r.Style = $origStyle
$ws.Range('E39').Value = '  +3.38%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$r = $ws.Range('D40')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.102'
$r.Style = $origStyle
$ws.Range('E40').Value = '  -2.01%  '
$r = $ws.Range('D41')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '6.42'
$r.Style = $origStyle
$ws.Range('E41').Value = '  +5.99%  '
$ws.Range('E42').Value = '  -1.44%  '
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('E44').Value = '  -0.01%  '
$r = $ws.Range('D45')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '171.06'
$r.Style = $origStyle
$ws.Range('E45').Value = '  +4.63%  '
$ws.Range('E46').Value = '  +0.11%  '
$r = $ws.Range('D47')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '0.000311'
$r.Style = $origStyle
$ws.Range('E47').Value = '  -3.87%  '
$r = $ws.Range('D48')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '49.98'
$r.Style = $origStyle
$ws.Range('E48').Value = '  +2.65%  '
$r = $ws.Range('D49')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '428.63'
$r.Style = $origStyle
$ws.Range('E49').Value = '  +2.16%  '
$r = $ws.Range('D50')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '1.39'
$r.Style = $origStyle
$ws.Range('E50').Value = '  -0.46%  '
$r = $ws.Range('D51')
$origStyle = $r.Style
$r.NumberFormat = "@"
$r.Value = '8.64'
$r.Style = $origStyle
$ws.Range('E51').Value = '  +0.03%  '
